# replyBriefTemplate.docx - "fixed role on the word doc"
#
# The document has several "<Content Select="./CircuitCourtCase/Role" />"
# content-control placeholders used in the attorney signature blocks.
# Three of the four signature blocks already have a literal "-Appellant"
# suffix typed right after the content control (e.g. "...Role" />-Appellant"),
# but the very first signature block (right after the "Respectfully
# submitted," section) was missing that suffix. This fixes the
# inconsistency by adding the missing "-Appellant" run immediately after
# that content control, inside the same paragraph.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    # Identify the paragraph that holds only the bare
    # "<Content Select="./CircuitCourtCase/Role" />" placeholder (not the
    # "REPLY BRIEF OF <.../Role>" heading, and not the blocks that already
    # carry the "-Appellant" suffix).
    if ($t -like "*CircuitCourtCase/Role*" -and
        $t -notlike "*REPLY*" -and
        $t -notlike "*Appellant*") {

        # Collapse a range to the end of the paragraph (just before the
        # paragraph mark) and insert the missing text there, right after
        # the content control's closing tag.
        $r = $p.Range.Duplicate
        $r.Collapse(0)
        [void]$r.MoveEnd(1, -1)
        $r.Text = "-Appellant"
    }
}
